$wb = $excel.ActiveWorkbook

# Sheets (order in workbook): 1=Contact, 2=QuickLink, 3=ContactTypes, 4=Users
$wsContact = $wb.Worksheets.Item(1)
$wsQuickLink = $wb.Worksheets.Item(2)

# --- Contact sheet: update Quick Link column header values (text unchanged, kept for completeness) ---
$wsContact.Range("F1").Value = "Contact Full Name"
$wsContact.Range("F2").Value = "Test ExternalContact"

# --- QuickLink sheet: refresh the list of quick links ---
$wsQuickLink.Range("A1").Value = "Quick Link"
$wsQuickLink.Range("A2").Value = "HL Relationships"
$wsQuickLink.Range("A3").Value = "Industry Focus"
$wsQuickLink.Range("A4").Value = "Opportunity Contacts"
$wsQuickLink.Range("A5").Value = "Engagement Contacts (Contact)"
$wsQuickLink.Range("A6").Value = "Engagements Shown"
$wsQuickLink.Range("A7").Value = "Affiliated Companies"
$wsQuickLink.Range("A8").Value = "Related Companies"
$wsQuickLink.Range("A9").Value = "Memberships"
$wsQuickLink.Range("A10").Value = "Contact Sectors"
$wsQuickLink.Range("A11").Value = "Campaign History"
$wsQuickLink.Range("A12").Value = "Contact Email History"
$wsQuickLink.Range("A13").Value = "Contact Sources"
$wsQuickLink.Range("A14").Value = "Development Leads"
$wsQuickLink.Range("A15").Value = "Files"
$wsQuickLink.Range("A16").Value = "Contact History"

# --- Make QuickLink the active sheet/tab, with a new selection ---
$wsQuickLink.Activate()
$wsQuickLink.Range("A17").Select()
